# Add "Romansh" (rm) language support to the properties template.
#
# The header row currently reads:
#   name | super | object | en | de | fr | it |
#   comment_en | comment_de | comment_fr | comment_it |
#   gui_element | gui_attributes
#
# Two new columns are inserted:
#   - "rm"          right after "it"            (new col H)
#   - "comment_rm"  right after "comment_it"     (new col M, after the "rm" shift)
#
# Inserting a full column shifts everything to its right one column over and
# carries along the formatting of the columns being pushed right (matching
# Excel's normal "Insert" behaviour), which is exactly what the target sheet
# needs: the previously-filled data cells in rows 2 and 3 slide right with
# their styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column H for "rm" (pushes old H..M -> I..N)
$ws.Columns("H:H").Insert()
$ws.Range("H1").Value = "rm"

# 2) Insert a new column M for "comment_rm" (pushes the shifted L(comment_it).. -> onward)
$ws.Columns("M:M").Insert()
$ws.Range("M1").Value = "comment_rm"
